# Update countries & provincias Spain
# Refreshes the COVID-19 "paises" stats table: new timestamp in A1, updated
# totals for several countries, and a handful of countries whose row order
# swapped (so the corresponding label + stats move together) as their case
# counts changed rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Julio de 2020 a las 02:43"

# Row 4
$ws.Cells.Item(4, 2).Value = 4028362
$ws.Cells.Item(4, 3).Value = 66933
$ws.Cells.Item(4, 4).Value = 1885527
$ws.Cells.Item(4, 5).Value = 1997891
$ws.Cells.Item(4, 7).Value = 1110
$ws.Cells.Item(4, 8).Value = 144944

# Row 6
$ws.Cells.Item(6, 5).Value = 412922
$ws.Cells.Item(6, 7).Value = 671
$ws.Cells.Item(6, 8).Value = 28770

# Row 24
$ws.Cells.Item(24, 2).Value = 111697
$ws.Cells.Item(24, 3).Value = 573
$ws.Cells.Item(24, 4).Value = 97757
$ws.Cells.Item(24, 5).Value = 5078

# Row 41
$ws.Cells.Item(41, 1).Value = "Panama"
$ws.Cells.Item(41, 2).Value = 55153
$ws.Cells.Item(41, 3).Value = 727
$ws.Cells.Item(41, 4).Value = 30075
$ws.Cells.Item(41, 5).Value = 23919
$ws.Cells.Item(41, 7).Value = 32
$ws.Cells.Item(41, 8).Value = 1159

# Row 42
$ws.Cells.Item(42, 1).Value = "Republica Dominicana"
$ws.Cells.Item(42, 2).Value = 54797
$ws.Cells.Item(42, 3).Value = 841
$ws.Cells.Item(42, 4).Value = 25976
$ws.Cells.Item(42, 5).Value = 27822
$ws.Cells.Item(42, 7).Value = 6
$ws.Cells.Item(42, 8).Value = 999

# Row 70
$ws.Cells.Item(70, 1).Value = "Chequia"
$ws.Cells.Item(70, 2).Value = 14324
$ws.Cells.Item(70, 3).Value = 226
$ws.Cells.Item(70, 4).Value = 8918
$ws.Cells.Item(70, 5).Value = 5046
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 360

# Row 71
$ws.Cells.Item(71, 1).Value = "Kenia"
$ws.Cells.Item(71, 2).Value = 14168
$ws.Cells.Item(71, 3).Value = 397
$ws.Cells.Item(71, 4).Value = 6258
$ws.Cells.Item(71, 5).Value = 7660
$ws.Cells.Item(71, 7).Value = 12
$ws.Cells.Item(71, 8).Value = 250

# Row 107
$ws.Cells.Item(107, 2).Value = 3135
$ws.Cells.Item(107, 3).Value = 5
$ws.Cells.Item(107, 4).Value = 1464
$ws.Cells.Item(107, 5).Value = 1578

# Row 113
$ws.Cells.Item(113, 2).Value = 2824
$ws.Cells.Item(113, 3).Value = 16
$ws.Cells.Item(113, 4).Value = 2614
$ws.Cells.Item(113, 5).Value = 172

# Row 117
$ws.Cells.Item(117, 2).Value = 2381
$ws.Cells.Item(117, 3).Value = 98
$ws.Cells.Item(117, 4).Value = 426
$ws.Cells.Item(117, 5).Value = 1920

# Row 129
$ws.Cells.Item(129, 1).Value = "Zimbabue"
$ws.Cells.Item(129, 2).Value = 1820
$ws.Cells.Item(129, 3).Value = 107
$ws.Cells.Item(129, 4).Value = 488
$ws.Cells.Item(129, 5).Value = 1306
$ws.Cells.Item(129, 8).Value = 26

# Row 130
$ws.Cells.Item(130, 1).Value = "Sierra Leona"
$ws.Cells.Item(130, 2).Value = 1727
$ws.Cells.Item(130, 3).Value = 16
$ws.Cells.Item(130, 4).Value = 1273
$ws.Cells.Item(130, 5).Value = 388
$ws.Cells.Item(130, 8).Value = 66

# Row 139
$ws.Cells.Item(139, 1).Value = "Surinam"
$ws.Cells.Item(139, 2).Value = 1131
$ws.Cells.Item(139, 3).Value = 52
$ws.Cells.Item(139, 4).Value = 705
$ws.Cells.Item(139, 5).Value = 405
$ws.Cells.Item(139, 8).Value = 21

# Row 140
$ws.Cells.Item(140, 1).Value = "Jordania"
$ws.Cells.Item(140, 2).Value = 1113
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 1034
$ws.Cells.Item(140, 5).Value = 68
$ws.Cells.Item(140, 8).Value = 11

# Row 141
$ws.Cells.Item(141, 2).Value = 1113
$ws.Cells.Item(141, 3).Value = 8
$ws.Cells.Item(141, 4).Value = 1018
$ws.Cells.Item(141, 5).Value = 26

# Row 142
$ws.Cells.Item(142, 1).Value = "Liberia"
$ws.Cells.Item(142, 2).Value = 1108
$ws.Cells.Item(142, 3).Value = 1
$ws.Cells.Item(142, 4).Value = 575
$ws.Cells.Item(142, 5).Value = 463
$ws.Cells.Item(142, 8).Value = 70

# Row 143
$ws.Cells.Item(143, 1).Value = "Uruguay"
$ws.Cells.Item(143, 2).Value = 1096
$ws.Cells.Item(143, 3).Value = 32
$ws.Cells.Item(143, 4).Value = 929
$ws.Cells.Item(143, 5).Value = 134
$ws.Cells.Item(143, 8).Value = 33

# Row 144
$ws.Cells.Item(144, 1).Value = "Uganda"
$ws.Cells.Item(144, 2).Value = 1072
$ws.Cells.Item(144, 3).Value = 3
$ws.Cells.Item(144, 4).Value = 958
$ws.Cells.Item(144, 5).Value = 114
$ws.Cells.Item(144, 8).Value = 0

# Row 145
$ws.Cells.Item(145, 1).Value = "Burkina Faso"
$ws.Cells.Item(145, 2).Value = 1065
$ws.Cells.Item(145, 4).Value = 901
$ws.Cells.Item(145, 5).Value = 111
$ws.Cells.Item(145, 8).Value = 53

# Row 153
$ws.Cells.Item(153, 4).Value = 588
$ws.Cells.Item(153, 5).Value = 144

# Row 162
$ws.Cells.Item(162, 2).Value = 401
$ws.Cells.Item(162, 3).Value = 17
$ws.Cells.Item(162, 5).Value = 36

# Row 166
$ws.Cells.Item(166, 2).Value = 339
$ws.Cells.Item(166, 3).Value = 2
$ws.Cells.Item(166, 5).Value = 157

# Row 169
$ws.Cells.Item(169, 2).Value = 328
$ws.Cells.Item(169, 3).Value = 6
$ws.Cells.Item(169, 5).Value = 120

# Row 173
$ws.Cells.Item(173, 4).Value = 201
$ws.Cells.Item(173, 5).Value = 1

# Row 210
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"

# Row 211
$ws.Cells.Item(211, 1).Value = "Groenlandia"
